# Replace the "180R" comment (resistor R3) with "470R".
# Row 8 holds designator "R3" in column A, with its value/comment in column C.
# A leading apostrophe is used so the cell keeps its existing "stored as text"
# (quote-prefix) formatting instead of Excel re-evaluating it as a plain value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "'470R"
